$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values for rows 2-252 based on the new run results.
$ws.Range("C2:C4").Value = 8147
$ws.Range("C5:C26").Value = 8056
$ws.Range("C27:C35").Value = 7619
$ws.Range("C36:C252").Value = 7586
